# Weekly update: insert a new "Orégano" price record for the week of
# 2022-04-20 (serial 44671) right after the current latest entries,
# pushing the existing rows 27..60 down to 28..61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 27 (shifts old rows 27-60 -> 28-61,
# carrying their formatting, e.g. the date-format style on column D).
$ws.Rows.Item(27).Insert()

# Populate the new row 27 with this week's record.
$ws.Cells.Item(27, 1).Value  = 9
$ws.Cells.Item(27, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(27, 3).Value  = "Metropolitana"
$ws.Cells.Item(27, 4).Value  = 44671
$ws.Cells.Item(27, 5).Value  = 13
$ws.Cells.Item(27, 6).Value  = 100112029
$ws.Cells.Item(27, 7).Value  = "Orégano"
$ws.Cells.Item(27, 8).Value  = "Sin especificar"
$ws.Cells.Item(27, 9).Value  = "Primera"
$ws.Cells.Item(27, 10).Value = 16
$ws.Cells.Item(27, 11).Value = 16000
$ws.Cells.Item(27, 12).Value = 16000
$ws.Cells.Item(27, 13).Value = 16000
$ws.Cells.Item(27, 14).Value = "$/docena de atados"
$ws.Cells.Item(27, 15).Value = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value = 5333
$ws.Cells.Item(27, 17).Value = 3
$ws.Cells.Item(27, 18).Value = "Hortaliza"
